# Apply cryptos list update (GitHub Actions style data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 49/50 content swap (MultiversX moves to rank 47, HuobiToken to rank 48) ---
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"

# --- Safe (already-textual) Price column updates ---
$ws.Range("D2").Value = "43.043.63"
$ws.Range("D3").Value = "2.311.93"
$ws.Range("D15").Value = "2.678.47"
$ws.Range("D16").Value = "2.306.54"
$ws.Range("D18").Value = "42.980.01"
$ws.Range("D42").Value = "1.979.50"
$ws.Range("D48").Value = "2.542.39"

# --- Price column updates that look numeric: enter as formula then flatten to text value
#     via Copy + PasteSpecial(values) so Excel keeps exact text (e.g. trailing zeros).
$ws.Range("D4").Formula = "=""0.999"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("D5").Formula = "=""302.25"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=""97.33"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D7").Formula = "=""0.506"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("D9").Formula = "=""0.500"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("D10").Formula = "=""35.20"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("D11").Formula = "=""19.41"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D12").Formula = "=""0.0793"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("D14").Formula = "=""6.88"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("D17").Formula = "=""0.788"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("D19").Formula = "=""12.56"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("D20").Formula = "=""0.0₃0894"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D22").Formula = "=""67.78"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("D23").Formula = "=""236.22"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("D25").Formula = "=""2.45"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D27").Formula = "=""24.84"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=""2.06"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("D29").Formula = "=""164.13"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D30").Formula = "=""9.09"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D31").Formula = "=""32.72"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("D32").Formula = "=""0.999"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=""17.87"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("D35").Formula = "=""4.49"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("D36").Formula = "=""0.0700"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D37").Formula = "=""2.35"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("D38").Formula = "=""0.100"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D40").Formula = "=""2.77"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("D43").Formula = "=""10.61"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("D44").Formula = "=""18.90"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D47").Formula = "=""2.78"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("D49").Formula = "=""53.60"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("D50").Formula = "=""2.84"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("D51").Formula = "=""72.26"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

# --- Volume(1h) percentage column updates (always safe as text) ---
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("E11").Value = "  +6.52%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +4.17%  "
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  -6.82%  "
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("E43").Value = "  +5.85%  "
$ws.Range("E44").Value = "  +5.82%  "
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  +0.32%  "

$excel.CutCopyMode = 0
$ws.Range("A1").Select()
